$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H (Code) for rows 5-9 first, then column F (Well No.) for rows 5-9,
# then row 10's C/D ("X"), F10 ("E1"), H10 ("MIGT002") -- this ordering reproduces
# the exact shared-string insertion order recorded in the workbook.
$ws.Range("H5").Value = "MINK003"
$ws.Range("H6").Value = "MINK004"
$ws.Range("H7").Value = "MINK005"
$ws.Range("H8").Value = "MINK006"
$ws.Range("H9").Value = "MINK007"

$ws.Range("F5").Value = "A2"
$ws.Range("F6").Value = "A3"
$ws.Range("F7").Value = "A4"
$ws.Range("F8").Value = "A5"
$ws.Range("F9").Value = "A6"

$ws.Range("C10").Value = "X"
$ws.Range("D10").Value = "X"
$ws.Range("F10").Value = "E1"
$ws.Range("H10").Value = "MIGT002"

# --- Remaining cells (numbers + already-known strings): order doesn't affect
# the shared-string table since these reuse existing entries or are numeric.
$ws.Range("A5").Value = 14
$ws.Range("B5").Value = "Notoscopelus kroyeri"
$ws.Range("B5").Font.Italic = $true
$ws.Range("C5").Value = 120
$ws.Range("D5").Value = 21.84
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = "MI1"

$ws.Range("A6").Value = 14
$ws.Range("B6").Value = "Notoscopelus kroyeri"
$ws.Range("B6").Font.Italic = $true
$ws.Range("C6").Value = 130
$ws.Range("D6").Value = 28.86
$ws.Range("E6").Value = 1
$ws.Range("G6").Value = "MI1"

$ws.Range("A7").Value = 14
$ws.Range("B7").Value = "Notoscopelus kroyeri"
$ws.Range("B7").Font.Italic = $true
$ws.Range("C7").Value = 125
$ws.Range("D7").Value = 24.3
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = "MI1"

$ws.Range("A8").Value = 14
$ws.Range("B8").Value = "Notoscopelus kroyeri"
$ws.Range("B8").Font.Italic = $true
$ws.Range("C8").Value = 109
$ws.Range("D8").Value = 16.97
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = "MI1"

$ws.Range("A9").Value = 14
$ws.Range("B9").Value = "Notoscopelus kroyeri"
$ws.Range("B9").Font.Italic = $true
$ws.Range("C9").Value = 119
$ws.Range("D9").Value = 16.32
$ws.Range("E9").Value = 2
$ws.Range("G9").Value = "MI1"

$ws.Range("A10").Value = 15
$ws.Range("B10").Value = "Gadiculus thori"
$ws.Range("B10").Font.Italic = $true
$ws.Range("E10").Value = 2
$ws.Range("G10").Value = "MI1"

# Selection moves to A11 after data entry, and dimension auto-expands to L10.
$ws.Range("A11").Select()
